$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column N: header ---------------------------------------------
$ws.Range("N1").Value2 = "TS_TestDesciption"

# Style N1 like the other header cells (K1: bold Calibri on yellow fill,
# centered + wrapped) but without a border, so copy K1's format then
# strip the border back off.
$ws.Range("K1").Copy() | Out-Null
$ws.Range("N1").PasteSpecial(-4122) | Out-Null
$ws.Range("N1").Borders.LineStyle = -4142

# --- New column N: row 2 (long precondition/test description) --------
$n2Text = @"
**Precondition: 
	1. Flashing of “XCPACTDEV” file*
	(ex.
	a.	5G3: FL_3WA907541F_<Release>_XCPACTDEV_XXXX_E.pdx
	b.	5G5: FL_3WA907670F_<Release>_XCPACTDEV_XXXX_E.pdx)
	2. SFD: available
	3. Deactive all PDA group before check for each case
CASE 1:
Set: "WriteDataByIdentifierRequest[Group 3 Xcp].bit 5" == 'Enabled'
Check: the Sensor(Hella)Development Messages is sent on VCAN-Bus
CASE 2:
Set: "WriteDataByIdentifierRequest[Group 3 Xcp].bit 5" <> 'Enabled'
Check: Sending Sensor(Hella)Development Messages shall be stopped on VCAN.
"@

$ws.Range("N2").Value2 = $n2Text
$ws.Range("N2").WrapText = $true

# Wrapping that much text auto-grows the row; put row 2 back to its
# original height like the rest of the data rows.
$ws.Rows.Item(2).RowHeight = 15

# --- New column N: remaining rows get an (empty) cell, matching the ---
# --- placeholder empty cells already present in columns L/M -----------
$ws.Range("N3:N14").Borders.LineStyle = -4142

# --- Column widths (L widened slightly, N given its own width) -------
$ws.Columns.Item(12).ColumnWidth = 15.1
$ws.Columns.Item(14).ColumnWidth = 16.85

# --- Selection follows the freshly-entered cell, like the source edit -
$ws.Range("N2").Select() | Out-Null
